# Apply cell-value updates per the crypto price refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.283.59'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '1.790.87'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = "'226.03"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('E6').Value = '  +1.08%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = "'32.35"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('D10').Value = "'0.0689"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').Value = '2.048.64'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.818.91'
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = "'11.05"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.70%  '
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('D16').Value = '34.269.88'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').Value = "'68.06"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = '0.0₃0802'
$ws.Range('E19').Value = '  +2.34%  '
$ws.Range('D20').Value = "'246.84"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = "'10.95"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').Value = "'4.17"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('D25').Value = "'162.32"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').Value = "'16.37"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  +1.26%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = "'1.23"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = "'0.0521"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('D32').Value = "'3.76"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.14%  '
$ws.Range('D33').Value = "'3.85"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.89%  '
$ws.Range('D34').Value = "'1.81"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('D35').Value = '1.442.75'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('E36').Value = '  +8.37%  '
$ws.Range('D37').Value = "'0.661"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('D39').Value = "'0.0190"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('D40').Value = "'82.42"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.67%  '
$ws.Range('E41').Value = '  +2.06%  '
$ws.Range('D42').Value = "'14.14"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.53%  '
$ws.Range('D44').Value = "'0.923"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').Value = "'0.0520"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.84%  '
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('D48').Value = '1.944.21'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('D49').Value = "'105.49"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.45%  '
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E51').Value = '  -7.22%  '
